$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.378.13"
$ws.Range("E2").Value = "  +2.60%  "

$ws.Range("D3").Value = "1.822.87"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.81"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4683"
$ws.Range("E7").Value = "  +5.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3780"
$ws.Range("E8").Value = "  +3.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07445"
$ws.Range("E9").Value = "  +2.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8756"
$ws.Range("E10").Value = "  +2.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.76"
$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("D12").Value = "1.825.61"
$ws.Range("E12").Value = "  -15.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.689"
$ws.Range("E13").Value = "  +1.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.417"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.74"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07096"

$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008805"
$ws.Range("E18").Value = "  +1.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("E20").Value = "  +1.67%  "

$ws.Range("D21").Value = "27.398.40"
$ws.Range("E21").Value = "  +2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  +3.70%  "

$ws.Range("E23").Value = "  +2.05%  "

$ws.Range("D24").Value = "2.055.80"
$ws.Range("E24").Value = "  -8.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.936"
$ws.Range("E25").Value = "  -2.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.21"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.261"
$ws.Range("E27").Value = "  +4.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.60"
$ws.Range("E28").Value = "  +1.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.319"
$ws.Range("E29").Value = "  +2.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.41"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08927"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7844"
$ws.Range("E32").Value = "  +6.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.192"
$ws.Range("E33").Value = "  +3.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.533"
$ws.Range("E34").Value = "  +2.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.949"
$ws.Range("E35").Value = "  +1.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.100"
$ws.Range("E37").Value = "  +1.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01973"
$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05245"
$ws.Range("E39").Value = "  +1.70%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5355"
$ws.Range("E40").Value = "  +1.09%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.244"
$ws.Range("E41").Value = "  +3.46%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.385"
$ws.Range("E42").Value = "  +22.72%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.908"
$ws.Range("E43").Value = "  +2.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1699"
$ws.Range("E44").Value = "  +1.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.642"
$ws.Range("E45").Value = "  +3.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5084"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.54"
$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.61"
$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.679"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9989"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06372"
$ws.Range("E51").Value = "  +1.28%  "
